$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.623.33"
$ws.Range("E2").Value = "  +5.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.474.90"
$ws.Range("E3").Value = "  +3.90%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.58"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.94"
$ws.Range("E6").Value = "  +15.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.466.82"
$ws.Range("E7").Value = "  +3.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +2.13%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.694"
$ws.Range("E10").Value = "  +9.25%  "

$ws.Range("E11").Value = "  +29.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.74"
$ws.Range("E12").Value = "  +8.63%  "

$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.016.10"
$ws.Range("E14").Value = "  +3.65%  "

$ws.Range("E15").Value = "  +2.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.26"
$ws.Range("E16").Value = "  +4.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.527.84"
$ws.Range("E17").Value = "  +5.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.505.40"
$ws.Range("E18").Value = "  +5.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.05"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.15"
$ws.Range("E20").Value = "  +3.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000137"
$ws.Range("E21").Value = "  +21.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.38"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.23"
$ws.Range("E23").Value = "  +0.59%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.34"
$ws.Range("E24").Value = "  +8.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "312.54"
$ws.Range("E25").Value = "  +2.43%  "

$ws.Range("E26").Value = "  +0.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.42"
$ws.Range("E27").Value = "  +6.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.16"
$ws.Range("E28").Value = "  +3.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("E29").Value = "  +3.27%  "

$ws.Range("E30").Value = "  +2.78%  "

$ws.Range("B31").Value = "LEO"
$ws.Range("C31").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.38"
$ws.Range("E31").Value = "  -2.02%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  +3.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "44.89"
$ws.Range("E33").Value = "  +10.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.16"
$ws.Range("E34").Value = "  +4.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.67"
$ws.Range("E35").Value = "  +25.42%  "

$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0495"
$ws.Range("E37").Value = "  -5.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.54"
$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("E39").Value = "  +5.06%  "

$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  -1.73%  "

$ws.Range("E42").Value = "  +3.47%  "

$ws.Range("E43").Value = "  +2.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "137.82"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("E45").Value = "  +5.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.291"
$ws.Range("E46").Value = "  +4.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.02"
$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("E48").Value = "  +0.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.58"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.810.86"
$ws.Range("E51").Value = "  +3.53%  "
